$d = $word.ActiveDocument

$replacements = @(
    @{old="846×4=3384"; new="788×3=2364"},
    @{old="172×9=1548"; new="308×5=1540"},
    @{old="398×4=1592"; new="634×5=3170"},
    @{old="419×7=2933"; new="170×3=510"},
    @{old="613×7=4291"; new="137×6=822"},
    @{old="965×4=3860"; new="831×2=1662"},
    @{old="393×5=1965"; new="122×3=366"},
    @{old="179×9=1611"; new="304×3=912"},
    @{old="941×6=5646"; new="581×4=2324"},
    @{old="243×2=486";  new="991×6=5946"},
    @{old="494×4=1976"; new="928×2=1856"},
    @{old="750×8=6000"; new="346×7=2422"},
    @{old="749×4=2996"; new="402×5=2010"},
    @{old="621×5=3105"; new="259×3=777"},
    @{old="750×2=1500"; new="378×7=2646"},
    @{old="416×2=832";  new="840×3=2520"},
    @{old="242×4=968";  new="205×5=1025"},
    @{old="580×2=1160"; new="230×7=1610"},
    @{old="584×2=1168"; new="299×7=2093"},
    @{old="782×6=4692"; new="380×4=1520"},
    @{old="816×4=3264"; new="203×4=812"},
    @{old="131×2=262";  new="968×4=3872"},
    @{old="748×6=4488"; new="798×9=7182"},
    @{old="709×4=2836"; new="901×3=2703"},
    @{old="787×7=5509"; new="958×8=7664"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
